$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New log entry appended as row 49 (previous last data row was 48)
$newRow = 49

# Duplicate the formatting of the last existing data row (row 48) onto the new row
$ws.Range("A48:H48").Copy() | Out-Null
$ws.Range("A49:H49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = "2025-08-23 12:59:14 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-23 18:29:14 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
